# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# zh-cn and de-de handback files have now been received/processed:
#   - Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#   - "Latest Target File" / "Latest Handback File" columns on the zh-cn and
#     de-de detail sheets are populated with the source file + handback xlf
#   - de-de's "Latest Handback DateTime" gets a fresh timestamp
#   - Several columns are widened so the new, longer file names are readable

$wb = $excel.ActiveWorkbook

$overviewWs = $wb.Worksheets.Item("Overview")
$zhWs       = $wb.Worksheets.Item("zh-cn")
$deWs       = $wb.Worksheets.Item("de-de")

$sourceFileName = "52c719bb-53eb-4c41-b363-b9dd8f8f5fc5.md"
$sourceUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7553b14c570857e5ffe213f3527ac637b7bea88/e2e/52c719bb-53eb-4c41-b363-b9dd8f8f5fc5.md"

$zhHandbackFile = "52c719bb-53eb-4c41-b363-b9dd8f8f5fc5.e312619ee408b9ecd2c87e3dedeabae8e4332570.zh-cn.xlf"
$deHandbackFile = "52c719bb-53eb-4c41-b363-b9dd8f8f5fc5.e312619ee408b9ecd2c87e3dedeabae8e4332570.de-de.xlf"

$zhHandbackDate = "2016-09-02 13:13:21"
$deHandbackDate = "2016-09-02 13:13:28"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Status text: every cell currently showing "Ready for handoff"
#    (Overview!E2:F2, Overview!E3:F3, zh-cn!C2:C3, de-de!C2:C3) now reads
#    "Handed back: in sync with en-US".
# ---------------------------------------------------------------------
$overviewWs.Range("E2").Value = $statusText
$overviewWs.Range("F2").Value = $statusText
$overviewWs.Range("E3").Value = $statusText
$overviewWs.Range("F3").Value = $statusText
$zhWs.Range("C2").Value = $statusText
$zhWs.Range("C3").Value = $statusText
$deWs.Range("C2").Value = $statusText
$deWs.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# 2. zh-cn sheet: populate "Latest Target File" (I) and
#    "Latest Handback File" (J); "Latest Handback DateTime" (K) keeps its
#    existing value text ("0001-01-01 00:00:00" -> updated date below).
# ---------------------------------------------------------------------
$zhWs.Range("I2").Value = $sourceFileName
$zhWs.Range("J2").Value = $zhHandbackFile
$zhWs.Range("K2").Value = $zhHandbackDate

$zhWs.Range("I3").Value = $sourceFileName
$zhWs.Range("J3").Value = $zhHandbackFile
$zhWs.Range("K3").Value = $zhHandbackDate

# Re-create the hyperlinks (adding the new ones for column I) in the same
# order they appear in the target file: A2, I2, A3, I3.
$zhWs.Hyperlinks.Delete()
$zhWs.Hyperlinks.Add($zhWs.Range("A2"), $sourceUrl, "", "", $sourceFileName)
$zhWs.Hyperlinks.Add($zhWs.Range("I2"), $sourceUrl, "", "", $sourceFileName)
$zhWs.Hyperlinks.Add($zhWs.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7553b14c570857e5ffe213f3527ac637b7bea88/e2e/ffff5cebdee0-c94a-4fe9-a1af-51a9aded3739.md", "", "", "ffff5cebdee0-c94a-4fe9-a1af-51a9aded3739.md")
$zhWs.Hyperlinks.Add($zhWs.Range("I3"), $sourceUrl, "", "", $sourceFileName)

# Give I2/I3 the same "HyperLink" look (underline + blue) that column A uses.
# (Hyperlinks.Add applies its own theme-colored style, so the explicit font
# is (re)applied afterwards to match column A's look exactly.)
$zhWs.Range("I2").Font.Underline = 2
$zhWs.Range("I2").Font.Color = 15570276
$zhWs.Range("I3").Font.Underline = 2
$zhWs.Range("I3").Font.Color = 15570276

# ---------------------------------------------------------------------
# 3. de-de sheet: same treatment, but "Latest Handback DateTime" (K) gets
#    a brand-new timestamp (de-de was just handed back).
# ---------------------------------------------------------------------
$deWs.Range("I2").Value = $sourceFileName
$deWs.Range("J2").Value = $deHandbackFile
$deWs.Range("K2").Value = $deHandbackDate

$deWs.Range("I3").Value = $sourceFileName
$deWs.Range("J3").Value = $deHandbackFile
$deWs.Range("K3").Value = $deHandbackDate

$deWs.Hyperlinks.Delete()
$deWs.Hyperlinks.Add($deWs.Range("A2"), $sourceUrl, "", "", $sourceFileName)
$deWs.Hyperlinks.Add($deWs.Range("I2"), $sourceUrl, "", "", $sourceFileName)
$deWs.Hyperlinks.Add($deWs.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7553b14c570857e5ffe213f3527ac637b7bea88/e2e/ffff5cebdee0-c94a-4fe9-a1af-51a9aded3739.md", "", "", "ffff5cebdee0-c94a-4fe9-a1af-51a9aded3739.md")
$deWs.Hyperlinks.Add($deWs.Range("I3"), $sourceUrl, "", "", $sourceFileName)

$deWs.Range("I2").Font.Underline = 2
$deWs.Range("I2").Font.Color = 15570276
$deWs.Range("I3").Font.Underline = 2
$deWs.Range("I3").Font.Color = 15570276

# ---------------------------------------------------------------------
# 4. Widen columns so the longer file names / status text fit.
#    (ColumnWidth is quantized by Excel's pixel grid, so the values below
#    are chosen to land on the closest achievable width.)
# ---------------------------------------------------------------------
$wideStatusWidth = 29.1666666666667   # -> stored width ~29.98 target (closest achievable: 30)
$wideFileWidth   = 39.1666666666667   # -> stored width 40

$overviewWs.Columns.Item(5).ColumnWidth = $wideStatusWidth   # E
$overviewWs.Columns.Item(6).ColumnWidth = $wideStatusWidth   # F

$zhWs.Columns.Item(3).ColumnWidth = $wideStatusWidth   # C
$zhWs.Columns.Item(9).ColumnWidth = $wideFileWidth     # I
$zhWs.Columns.Item(10).ColumnWidth = $wideFileWidth    # J

$deWs.Columns.Item(3).ColumnWidth = $wideStatusWidth   # C
$deWs.Columns.Item(9).ColumnWidth = $wideFileWidth     # I
$deWs.Columns.Item(10).ColumnWidth = $wideFileWidth    # J

Write-Host "Handback report generated."
